$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Status column: rows 4 and 5 (PageHeader, ToastMessageValidation)
# change from "FAIL" to "PASS"
$ws.Range("D4").Value = "PASS"
$ws.Range("D5").Value = "PASS"
